# AF-611: tests for named ranges and area intersections are moved to
# temporary_excel_files.
#
# The ISEVEN test in D12 no longer exercises the named range "hjk"
# (Sheet1!$D$3) - it now mirrors the plain-cell-reference case in D7,
# and the active selection moves on to the next (empty) row beneath
# the test table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D12 used to test ISEVEN against the defined name "hjk"; repoint it at
# a plain cell reference instead (A1 is odd, so the result flips to FALSE).
$ws.Range("D12").Formula = "=ISEVEN(A1)"

# Move the active cell/selection down to D13, just past the last test row.
$ws.Range("D13").Select()
